# Scheduled-runner update: refresh computed leve-profit columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-class sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 427.18182
$ws.Range("I8").Value = 44.75
$ws.Range("J8").Value = 645.7143
$ws.Range("K8").Value = 134.25
$ws.Range("L8").Value = 1937.1429
$ws.Range("M8").Value = 4.75
$ws.Range("N8").Value = -2215.1429
$ws.Range("H33").Value = 278.62857
$ws.Range("I33").Value = 239.03334
$ws.Range("K33").Value = 239.03334
$ws.Range("M33").Value = -10.03334000000001
$ws.Range("H43").Value = 599
$ws.Range("J43").Value = 599
$ws.Range("L43").Value = 599
$ws.Range("N43").Value = -737
$ws.Range("H58").Value = 592.6316
$ws.Range("J58").Value = 9990
$ws.Range("L58").Value = 29970
$ws.Range("N58").Value = -30270
$ws.Range("H62").Value = 2501
$ws.Range("I62").Value = 2476.25
$ws.Range("K62").Value = 2476.25
$ws.Range("M62").Value = -1852.25
$ws.Range("H65").Value = 2501
$ws.Range("I65").Value = 2476.25
$ws.Range("K65").Value = 12381.25
$ws.Range("M65").Value = -9261.25
$ws.Range("H76").Value = 3320
$ws.Range("I76").Value = 3200
$ws.Range("K76").Value = 3200
$ws.Range("M76").Value = -2885
$ws.Range("H79").Value = 3320
$ws.Range("I79").Value = 3200
$ws.Range("K79").Value = 3200
$ws.Range("M79").Value = -2108
$ws.Range("H88").Value = 1732.2222
$ws.Range("I88").Value = 1101.5
$ws.Range("J88").Value = 1912.4286
$ws.Range("K88").Value = 1101.5
$ws.Range("L88").Value = 1912.4286
$ws.Range("M88").Value = -695.5
$ws.Range("N88").Value = -2724.4286
$ws.Range("H91").Value = 1732.2222
$ws.Range("I91").Value = 1101.5
$ws.Range("J91").Value = 1912.4286
$ws.Range("K91").Value = 1101.5
$ws.Range("L91").Value = 1912.4286
$ws.Range("M91").Value = 302.5
$ws.Range("N91").Value = -4720.4286
$ws.Range("H106").Value = 1541.8182
$ws.Range("I106").Value = 1345.5
$ws.Range("K106").Value = 1345.5
$ws.Range("M106").Value = -714.5
$ws.Range("H116").Value = 2846.2222
$ws.Range("J116").Value = 3187.2856
$ws.Range("L116").Value = 3187.2856
$ws.Range("N116").Value = -10071.2856
$ws.Range("H129").Value = 879.6739
$ws.Range("J129").Value = 885.63635
$ws.Range("L129").Value = 2656.90905
$ws.Range("N129").Value = -12656.90905
$ws.Range("H132").Value = 37207.793
$ws.Range("I132").Value = 41447.152
$ws.Range("J132").Value = 466.66666
$ws.Range("K132").Value = 124341.456
$ws.Range("L132").Value = 1399.99998
$ws.Range("M132").Value = -121811.456
$ws.Range("N132").Value = -6459.999980000001
$ws.Range("H138").Value = 3173.7874
$ws.Range("J138").Value = 3620.5
$ws.Range("L138").Value = 10861.5
$ws.Range("N138").Value = -21141.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19519.508
$ws.Range("I32").Value = 21562.76
$ws.Range("K32").Value = 21562.76
$ws.Range("M32").Value = -21275.76
$ws.Range("H132").Value = 26092.38
$ws.Range("I132").Value = 1854.8182
$ws.Range("J132").Value = 52753.7
$ws.Range("K132").Value = 5564.4546
$ws.Range("L132").Value = 158261.1
$ws.Range("M132").Value = -3034.4546
$ws.Range("N132").Value = -163321.1

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3848.2222
$ws.Range("J94").Value = 5393.5454
$ws.Range("L94").Value = 5393.5454
$ws.Range("N94").Value = -6295.5454
$ws.Range("H99").Value = 904.61536
$ws.Range("I99").Value = 896.6667
$ws.Range("K99").Value = 896.6667
$ws.Range("M99").Value = 601.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7637.2266
$ws.Range("I31").Value = 14322.958
$ws.Range("J31").Value = 2104.2068
$ws.Range("K31").Value = 14322.958
$ws.Range("L31").Value = 2104.2068
$ws.Range("M31").Value = -14027.958
$ws.Range("N31").Value = -2694.2068
$ws.Range("H34").Value = 7637.2266
$ws.Range("I34").Value = 14322.958
$ws.Range("J34").Value = 2104.2068
$ws.Range("K34").Value = 14322.958
$ws.Range("L34").Value = 2104.2068
$ws.Range("M34").Value = -14120.958
$ws.Range("N34").Value = -2508.2068
$ws.Range("H107").Value = 682.25
$ws.Range("I107").Value = 682.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 682.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1237.75
$ws.Range("N107").Value = ""
$ws.Range("H134").Value = 841.0952
$ws.Range("I134").Value = 721
$ws.Range("K134").Value = 2163
$ws.Range("M134").Value = 372

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 638.0244
$ws.Range("I5").Value = 497.21054
$ws.Range("J5").Value = 759.63635
$ws.Range("K5").Value = 1491.63162
$ws.Range("L5").Value = 2278.90905
$ws.Range("M5").Value = -1379.63162
$ws.Range("N5").Value = -2502.90905
$ws.Range("H68").Value = 1217.7084
$ws.Range("J68").Value = 1418
$ws.Range("L68").Value = 4254
$ws.Range("N68").Value = -5876
$ws.Range("H71").Value = 1217.7084
$ws.Range("J71").Value = 1418
$ws.Range("L71").Value = 12762
$ws.Range("N71").Value = -20874
$ws.Range("H86").Value = 849
$ws.Range("J86").Value = 805
$ws.Range("L86").Value = 2415
$ws.Range("N86").Value = -4787
$ws.Range("H89").Value = 849
$ws.Range("J89").Value = 805
$ws.Range("L89").Value = 7245
$ws.Range("N89").Value = -19101
$ws.Range("H131").Value = 145754.98
$ws.Range("J131").Value = 154662.06
$ws.Range("L131").Value = 463986.18
$ws.Range("N131").Value = -474066.18
$ws.Range("H135").Value = 638.0244
$ws.Range("I135").Value = 497.21054
$ws.Range("J135").Value = 759.63635
$ws.Range("K135").Value = 4474.894859999999
$ws.Range("L135").Value = 6836.72715
$ws.Range("M135").Value = -1939.894859999999
$ws.Range("N135").Value = -11906.72715
$ws.Range("H139").Value = 14070.25
$ws.Range("I139").Value = 15242.429
$ws.Range("J139").Value = 5865
$ws.Range("K139").Value = 45727.287
$ws.Range("L139").Value = 17595
$ws.Range("M139").Value = -40587.287
$ws.Range("N139").Value = -27875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4695.8
$ws.Range("I70").Value = 4250
$ws.Range("J70").Value = 4993
$ws.Range("K70").Value = 4250
$ws.Range("L70").Value = 4993
$ws.Range("M70").Value = -3980
$ws.Range("N70").Value = -5533
$ws.Range("H73").Value = 4695.8
$ws.Range("I73").Value = 4250
$ws.Range("J73").Value = 4993
$ws.Range("K73").Value = 4250
$ws.Range("L73").Value = 4993
$ws.Range("M73").Value = -3314
$ws.Range("N73").Value = -6865

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 39470
$ws.Range("J108").Value = 39470
$ws.Range("L108").Value = 39470
$ws.Range("N108").Value = -47150
$ws.Range("H122").Value = 1468.1818
$ws.Range("I122").Value = 1342.4667
$ws.Range("K122").Value = 4027.4001
$ws.Range("M122").Value = -1577.4001
